$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 733, shifting old rows 733:774 down to 734:775
# (and bumping the sheet's used range from A1:D774 to A1:D775).
$ws.Rows.Item(733).Insert()

# Populate the newly inserted row with the new data point. The date column
# is plain text in this sheet (e.g. "2026/12/29"), so prefix with an
# apostrophe to stop it being auto-converted to a date value, then reset the
# cell style so no stray number-format/quote-prefix styling is left behind.
$ws.Cells.Item(733, 1).Value = "'2026/01/27"
$ws.Cells.Item(733, 1).Style = "Normal"
$ws.Cells.Item(733, 2).Value = "火"
$ws.Cells.Item(733, 3).Value = 8
$ws.Cells.Item(733, 4).Value = 197
